{"js": "const replacements = [\n  [\"2024-09-25 Wednesday\", \"2024-09-26 Thursday\"],\n  [\"252\u00f74=\", \"632\u00f74=\"],\n  [\"710\u00f74=\", \"184\u00f77=\"],\n  [\"662\u00f73=\", \"817\u00f77=\"],\n  [\"355\u00f72=\", \"958\u00f76=\"],\n  [\"303\u00f74=\", \"912\u00f76=\"],\n  [\"581\u00f74=\", \"481\u00f78=\"],\n  [\"395\u00f75=\", \"838\u00f72=\"],\n  [\"770\u00f72=\", \"262\u00f75=\"],\n  [\"335\u00f74=\", \"715\u00f76=\"],\n  [\"842\u00f79=\", \"595\u00f75=\"],\n  [\"592\u00f77=\", \"457\u00f76=\"],\n  [\"925\u00f77=\", \"966\u00f79=\"],\n  [\"157\u00f75=\", \"728\u00f73=\"],\n  [\"841\u00f72=\", \"555\u00f73=\"],\n  [\"464\u00f79=\", \"637\u00f73=\"],\n  [\"185\u00f72=\", \"993\u00f73=\"],\n  [\"180\u00f79=\", \"974\u00f79=\"],\n  [\"977\u00f75=\", \"389\u00f74=\"],\n  [\"239\u00f78=\", \"494\u00f76=\"],\n  [\"653\u00f76=\", \"668\u00f78=\"],\n  [\"294\u00f74=\", \"451\u00f79=\"],\n  [\"136\u00f72=\", \"871\u00f79=\"],\n  [\"251\u00f79=\", \"708\u00f78=\"],\n  [\"896\u00f78=\", \"268\u00f75=\"],\n  [\"152\u00f79=\", \"212\u00f74=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-09-25 Wednesday\"; New = \"2024-09-26 Thursday\" },\n    @{ Old = \"252\u00f74=\"; New = \"632\u00f74=\" },\n    @{ Old = \"710\u00f74=\"; New = \"184\u00f77=\" },\n    @{ Old = \"662\u00f73=\"; New = \"817\u00f77=\" },\n    @{ Old = \"355\u00f72=\"; New = \"958\u00f76=\" },\n    @{ Old = \"303\u00f74=\"; New = \"912\u00f76=\" },\n    @{ Old = \"581\u00f74=\"; New = \"481\u00f78=\" },\n    @{ Old = \"395\u00f75=\"; New = \"838\u00f72=\" },\n    @{ Old = \"770\u00f72=\"; New = \"262\u00f75=\" },\n    @{ Old = \"335\u00f74=\"; New = \"715\u00f76=\" },\n    @{ Old = \"842\u00f79=\"; New = \"595\u00f75=\" },\n    @{ Old = \"592\u00f77=\"; New = \"457\u00f76=\" },\n    @{ Old = \"925\u00f77=\"; New = \"966\u00f79=\" },\n    @{ Old = \"157\u00f75=\"; New = \"728\u00f73=\" },\n    @{ Old = \"841\u00f72=\"; New = \"555\u00f73=\" },\n    @{ Old = \"464\u00f79=\"; New = \"637\u00f73=\" },\n    @{ Old = \"185\u00f72=\"; New = \"993\u00f73=\" },\n    @{ Old = \"180\u00f79=\"; New = \"974\u00f79=\" },\n    @{ Old = \"977\u00f75=\"; New = \"389\u00f74=\" },\n    @{ Old = \"239\u00f78=\"; New = \"494\u00f76=\" },\n    @{ Old = \"653\u00f76=\"; New = \"668\u00f78=\" },\n    @{ Old = \"294\u00f74=\"; New = \"451\u00f79=\" },\n    @{ Old = \"136\u00f72=\"; New = \"871\u00f79=\" },\n    @{ Old = \"251\u00f79=\"; New = \"708\u00f78=\" },\n    @{ Old = \"896\u00f78=\"; New = \"268\u00f75=\" },\n    @{ Old = \"152\u00f79=\"; New = \"212\u00f74=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute(\n        $pair.Old,    # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $pair.New,    # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n}\n"}
